$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.21"
$ws.Range("G2").Value = "'7"
$ws.Range("D3").Value = "'23.02"
$ws.Range("G3").Value = "'7"
$ws.Range("D4").Value = "'6.266"
$ws.Range("G4").Value = "'7"
$ws.Range("D5").Value = "'0.06266"
$ws.Range("G5").Value = "'7"
$ws.Range("D6").Value = "'3.663"
$ws.Range("G6").Value = "'7"
$ws.Range("D7").Value = "'6.667"
$ws.Range("G7").Value = "'7"
$ws.Range("D8").Value = "'1.362"
$ws.Range("G8").Value = "'7"
$ws.Range("D9").Value = "'0.8301"
$ws.Range("G9").Value = "'7"
$ws.Range("D10").Value = "'0.01378"
$ws.Range("G10").Value = "'7"
$ws.Range("D11").Value = "'0.1626"
$ws.Range("G11").Value = "'7"
$ws.Range("D12").Value = "'0.08306"
$ws.Range("G12").Value = "'7"
$ws.Range("D13").Value = "'0.03441"
$ws.Range("G13").Value = "'7"
$ws.Range("D14").Value = "'0.03091"
$ws.Range("G14").Value = "'7"
$ws.Range("D15").Value = "'0.09313"
$ws.Range("G15").Value = "'7"
$ws.Range("D16").Value = "'3.847"
$ws.Range("G16").Value = "'7"
$ws.Range("D17").Value = "'0.001640"
$ws.Range("G17").Value = "'7"
$ws.Range("D18").Value = "'0.04765"
$ws.Range("G18").Value = "'7"
$ws.Range("D19").Value = "'0.006405"
$ws.Range("G19").Value = "'7"
$ws.Range("D20").Value = "'0.005688"
$ws.Range("G20").Value = "'7"
$ws.Range("D21").Value = "'0.001093"
$ws.Range("G21").Value = "'7"
$ws.Range("D22").Value = "'0.0001550"
$ws.Range("G22").Value = "'7"
$ws.Range("D23").Value = "'3.715"
$ws.Range("G23").Value = "'7"
$ws.Range("G24").Value = "'7"
$ws.Range("G25").Value = "'7"
$ws.Range("G26").Value = "'7"
$ws.Range("G27").Value = "'7"
$ws.Range("G28").Value = "'7"
$ws.Range("G29").Value = "'7"
$ws.Range("G30").Value = "'7"
$ws.Range("G31").Value = "'7"
$ws.Range("G32").Value = "'7"
$ws.Range("G33").Value = "'7"
$ws.Range("G34").Value = "'7"
$ws.Range("G35").Value = "'7"
$ws.Range("G36").Value = "'7"
$ws.Range("G37").Value = "'7"
$ws.Range("G38").Value = "'7"
$ws.Range("G39").Value = "'7"
$ws.Range("D40").Value = "'0.04698"
$ws.Range("G40").Value = "'7"
$ws.Range("D41").Value = "'0.007061"
$ws.Range("G41").Value = "'7"
$ws.Range("D42").Value = "'0.1164"
$ws.Range("G42").Value = "'7"
$ws.Range("G43").Value = "'7"
$ws.Range("D44").Value = "'0.01219"
$ws.Range("G44").Value = "'7"
$ws.Range("G45").Value = "'7"
$ws.Range("B46").Value = "'Kangarootoken"
$ws.Range("C46").Value = "'https://coinranking.com/coin/zkVNkSGwZ3+kangarootoken-gar"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'45KangarootokenGAR"
$ws.Range("G46").Value = "'7"
$ws.Range("B47").Value = "'CoinbaseStockToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "'0.7698"
$ws.Range("E47").Value = "'46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("G47").Value = "'7"
$ws.Range("B48").Value = "'BOLO"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.04289"
$ws.Range("E48").Value = "'47BOLOBOLOBestin24h"
$ws.Range("G48").Value = "'7"
$ws.Range("B49").Value = "'CryptobidCoin"
$ws.Range("C49").Value = "'https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
$ws.Range("D49").Value = "'0.00002299"
$ws.Range("E49").Value = "'48CryptobidCoinCBC"
$ws.Range("G49").Value = "'7"
$ws.Range("B50").Value = "'SpecialPowerGold"
$ws.Range("C50").Value = "'https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg"
$ws.Range("D50").Value = "'0.01240"
$ws.Range("E50").Value = "'49SpecialPowerGoldSPG"
$ws.Range("G50").Value = "'7"
$ws.Range("B51").Value = "'DigiFinexToken"
$ws.Range("C51").Value = "'https://coinranking.com/coin/rY6dWXQL4+digifinextoken-dft"
$ws.Range("D51").Value = "'--"
$ws.Range("E51").Value = "'50DigiFinexTokenDFT"
$ws.Range("G51").Value = "'7"
